# Add a "Supplier" column (K) to the product import template.
# - New header "Supplier" in K1
# - Every data row (2-7) gets a supplier id of 2
# - Row 2 height nudged to 15 (matches the template after the edit)
# - Selection moved to the newly added K7 cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("K1").Value = "Supplier"

# Supplier id for every product row
$ws.Range("K2:K7").Value = 2

# Cosmetic row-height tweak that came along with the edit
$ws.Rows.Item(2).RowHeight = 15

# Leave the selection on the last cell of the new column, like the diff shows
$ws.Range("K7").Select()

Write-Output "Added Supplier column (K) with header + values, done."
